$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.219.20'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '1.867.79'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '0.7103'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').Value = '241.85'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '0.3114'
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').Value = '0.07661'
$ws.Range('E9').Value = '  -3.50%  '
$ws.Range('D10').Value = '24.70'
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('D11').Value = '0.08369'
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('D12').Value = '1.869.76'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = '5.229'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').Value = '0.7110'
$ws.Range('D15').Value = '91.27'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '29.223.78'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').Value = '5.946'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = '243.66'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').Value = '0.000007835'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').Value = '2.114.85'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '7.863'
$ws.Range('E23').Value = '  -1.81%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').Value = '0.1636'
$ws.Range('D26').Value = '163.06'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Value = '8.964'
$ws.Range('E27').Value = '  -0.95%  '
$ws.Range('D28').Value = '18.51'
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('D30').Value = '1.317'
$ws.Range('E30').Value = '  -3.06%  '
$ws.Range('D31').Value = '4.398'
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('D32').Value = '4.244'
$ws.Range('E32').Value = '  +3.09%  '
$ws.Range('D33').Value = '0.05136'
$ws.Range('E33').Value = '  -2.56%  '
$ws.Range('D34').Value = '0.7963'
$ws.Range('E34').Value = '  +9.40%  '
$ws.Range('D35').Value = '1.914'
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('D36').Value = '1.167'
$ws.Range('E36').Value = '  -2.70%  '
$ws.Range('D37').Value = '2.688'
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('D38').Value = '0.01857'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('D39').Value = '2.707'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '1.155.83'
$ws.Range('E40').Value = '  -5.58%  '
$ws.Range('D41').Value = '6.387'
$ws.Range('E41').Value = '  +4.28%  '
$ws.Range('D42').Value = '0.8957'
$ws.Range('E42').Value = '  -1.57%  '
$ws.Range('D43').Value = '73.27'
$ws.Range('E43').Value = '  -0.75%  '
$ws.Range('D44').Value = '0.9995'
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D46').Value = '2.011.69'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').Value = '0.5167'
$ws.Range('E47').Value = '  -2.24%  '
$ws.Range('D48').Value = '1.782'
$ws.Range('E48').Value = '  -1.09%  '
$ws.Range('D49').Value = '9.343'
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('D51').Value = '0.4295'
$ws.Range('E51').Value = '  -0.72%  '
